$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header in H1, copying the formatting from the
# neighboring header cell (G1) so it matches the other headers (bold,
# bordered, centered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Add the corresponding data value (0) in H2, unstyled like the other
# numeric data cells.
$ws.Range("H2").Value = 0
